# Settlements sheet: normalize the "size" column values to lowercase
# (Large/Medium/Small -> large/medium/small) without touching anything else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").Value   = "large"
$ws.Range("B5:B9").Value   = "medium"
$ws.Range("B10:B21").Value = "small"

# Leave the selection where the last edit was made, matching the saved view.
[void]$ws.Range("B10:B21").Select()
